$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 130
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = ""
$ws.Range("H64").Value = 4050
$ws.Range("H67").Value = 4050
$ws.Range("H74").Value = 4428.5713
$ws.Range("I74").Value = 4000
$ws.Range("K74").Value = 4000
$ws.Range("M74").Value = -3064
$ws.Range("H75").Value = 22661.334
$ws.Range("J75").Value = 22661.334
$ws.Range("L75").Value = 22661.334
$ws.Range("N75").Value = -24533.334
$ws.Range("H77").Value = 4428.5713
$ws.Range("I77").Value = 4000
$ws.Range("K77").Value = 20000
$ws.Range("M77").Value = -15320
$ws.Range("H78").Value = 22661.334
$ws.Range("J78").Value = 22661.334
$ws.Range("L78").Value = 67984.00199999999
$ws.Range("N78").Value = -77344.00199999999
$ws.Range("H96").Value = 1999.2222
$ws.Range("I96").Value = 1160.3334
$ws.Range("J96").Value = 2418.6667
$ws.Range("K96").Value = 3481.0002
$ws.Range("L96").Value = 7256.000100000001
$ws.Range("M96").Value = -2108.0002
$ws.Range("N96").Value = -10002.0001
$ws.Range("H132").Value = 1372.25
$ws.Range("I132").Value = 1372.25
$ws.Range("K132").Value = 4116.75
$ws.Range("M132").Value = -1586.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7833.6665
$ws.Range("I32").Value = 7833.6665
$ws.Range("K32").Value = 7833.6665
$ws.Range("M32").Value = -7546.6665
$ws.Range("H45").Value = 4399
$ws.Range("I45").Value = 4399
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 4399
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -4022
$ws.Range("N45").Value = ""
$ws.Range("H97").Value = 194.6
$ws.Range("I97").Value = 194.6
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 194.6
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 301.4
$ws.Range("N97").Value = ""
$ws.Range("H110").Value = 559.5
$ws.Range("I110").Value = 559.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 559.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1485.5
$ws.Range("N110").Value = ""

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1131.3636
$ws.Range("I99").Value = 1133.4
$ws.Range("K99").Value = 1133.4
$ws.Range("M99").Value = 364.5999999999999
$ws.Range("H105").Value = 2183.077
$ws.Range("I105").Value = 1676
$ws.Range("K105").Value = 1676
$ws.Range("M105").Value = 71

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6897.8
$ws.Range("I31").Value = 3287
$ws.Range("J31").Value = 8101.4
$ws.Range("K31").Value = 3287
$ws.Range("L31").Value = 8101.4
$ws.Range("M31").Value = -2992
$ws.Range("N31").Value = -8691.4
$ws.Range("H34").Value = 6897.8
$ws.Range("I34").Value = 3287
$ws.Range("J34").Value = 8101.4
$ws.Range("K34").Value = 3287
$ws.Range("L34").Value = 8101.4
$ws.Range("M34").Value = -3085
$ws.Range("N34").Value = -8505.4
$ws.Range("H60").Value = 4475
$ws.Range("I60").Value = 4475
$ws.Range("K60").Value = 4475
$ws.Range("M60").Value = -3964

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 400
$ws.Range("I118").Value = 400
$ws.Range("K118").Value = 1200
$ws.Range("M118").Value = 43
$ws.Range("H138").Value = 4817
$ws.Range("I138").Value = 634
$ws.Range("J138").Value = 9000
$ws.Range("K138").Value = 1902
$ws.Range("L138").Value = 27000
$ws.Range("M138").Value = 3238
$ws.Range("N138").Value = -37280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4731.25
$ws.Range("I102").Value = 3679.3333
$ws.Range("J102").Value = 7887
$ws.Range("K102").Value = 3679.3333
$ws.Range("L102").Value = 7887
$ws.Range("M102").Value = -2057.3333
$ws.Range("N102").Value = -11131
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4299.5
$ws.Range("I46").Value = 4299.5
$ws.Range("K46").Value = 4299.5
$ws.Range("M46").Value = -4111.5
$ws.Range("H61").Value = 2004
$ws.Range("I61").Value = 2004
$ws.Range("K61").Value = 2004
$ws.Range("M61").Value = -1802
$ws.Range("H68").Value = 2536
$ws.Range("I68").Value = 2536
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2536
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1787
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 2536
$ws.Range("I71").Value = 2536
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 12680
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -8936
$ws.Range("N71").Value = ""
$ws.Range("H93").Value = 3499.6667
$ws.Range("I93").Value = 3499.6667
$ws.Range("K93").Value = 3499.6667
$ws.Range("M93").Value = -2251.6667
$ws.Range("H113").Value = 2004
$ws.Range("I113").Value = 2004
$ws.Range("K113").Value = 2004
$ws.Range("M113").Value = 166
$ws.Range("H136").Value = 6401.6
$ws.Range("I136").Value = 6401.6
$ws.Range("K136").Value = 19204.8
$ws.Range("M136").Value = -16654.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 33167.332
$ws.Range("I4").Value = 40001
$ws.Range("J4").Value = 19500
$ws.Range("K4").Value = 40001
$ws.Range("L4").Value = 19500
$ws.Range("M4").Value = -39888
$ws.Range("N4").Value = -19726
$ws.Range("H126").Value = 1292.6666
$ws.Range("I126").Value = 993
$ws.Range("J126").Value = 1592.3334
$ws.Range("K126").Value = 2979
$ws.Range("L126").Value = 4777.0002
$ws.Range("M126").Value = -509
$ws.Range("N126").Value = -9717.0002
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470
